$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2199
$ws1.Range("F3").Value = 924
$ws1.Range("F4").Value = 1741
$ws1.Range("F5").Value = 397

# Sheet "全部类型" (All types) - aggregated view with the same events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 2199
$ws4.Range("F5").Value = 924
$ws4.Range("F6").Value = 1741
$ws4.Range("F7").Value = 397
